$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Build_Altar_1"
$ws.Range("C3").Value = "Build_Arena_1"
$ws.Range("C4").Value = "Build_Camp_1"
$ws.Range("C5").Value = "Build_Gold_Mine_1"
$ws.Range("C6").Value = "Build_Item_Hourse_1"
$ws.Range("C7").Value = "Build_League_1"
$ws.Range("C8").Value = "Build_Magic_Hourse_1"
$ws.Range("C9").Value = "Build_Tower_1"
$ws.Range("C10").Value = "Build_Town_1"
